$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was recorded for "Perejil" at Lo Valledor market.
# It is inserted as a brand-new row at position 407, pushing every
# subsequent row (407-522) down by one (408-523) while keeping all of
# their original values intact.
$ws.Rows(407).Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Range("A407").Value = 6
$ws.Range("B407").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C407").Value = "Metropolitana"
$ws.Range("D407").Value = 44736
$ws.Range("E407").Value = 13
$ws.Range("F407").Value = 100112044
$ws.Range("G407").Value = "Perejil"
$ws.Range("H407").Value = "Sin especificar"
$ws.Range("I407").Value = "Primera"
$ws.Range("J407").Value = 110
$ws.Range("K407").Value = 16000
$ws.Range("L407").Value = 17000
$ws.Range("M407").Value = 16455
$ws.Range("N407").Value = "`$/docena de atados"
$ws.Range("O407").Value = "Región Metropolitana"
$ws.Range("P407").Value = 5485
$ws.Range("Q407").Value = 3
$ws.Range("R407").Value = "Hortaliza"
